# Popravljen ispis kvara sa prioritetom
# Append a new "kvar" (fault) record as row 13 of the KVAROVI sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KVAROVI")

$targetRow = 13
$rng = $ws.Range("A" + $targetRow + ":F" + $targetRow)

# Force text storage so values like "2024-01-04" and "18" are kept as
# literal text (matching the existing IDKV/VRKV/IdEl columns) instead of
# being auto-converted to a date serial / number by Excel's input parser.
$rng.NumberFormat = "@"

$ws.Range("A" + $targetRow).Value = "20240104084543_01"
$ws.Range("B" + $targetRow).Value = "2024-01-04"
$ws.Range("C" + $targetRow).Value = "U popravci"
$ws.Range("D" + $targetRow).Value = "Kvar generatora"
$ws.Range("E" + $targetRow).Value = "Radni se"
$ws.Range("F" + $targetRow).Value = "18"

# Drop the explicit "Text" number format again so the new cells fall back
# to the workbook's default (general) cell style, same as every other row.
$rng.ClearFormats()
